$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.087075292901118928,
    -0.0099999998793904865,
    -0.0089999998805634362,
    0.28399755406391236,
    -0.0059999998851711922,
    -0.0059999998812685362,
    -0.019999999856175066,
    -0.019999999854625194,
    -0.0059999998773792029,
    -0.0059999998757334083,
    -0.0044999998781776185,
    -0.005999999875139661,
    0.0059815756191587965,
    -0.01199999986238609,
    -0.005999999872421391,
    -0.0059999998719511005,
    -0.005999999871326267,
    -0.008999999865896946,
    -0.0089999998810581516,
    -0.0089999998800198711,
    -0.0089999998798431236,
    -0.0089999998797027914,
    -0.013605373143122712,
    -0.090578451443487751,
    -0.041999999816268563,
    -0.0059999998810162936,
    -0.0059999998806388177,
    -0.0059999998793269782,
    -0.011999999867994049,
    -0.019999999853638428,
    -0.01499999986231515,
    0.018086556170155532,
    -0.0059999998781368191
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
